$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 5 new result rows under the existing header row (row 1),
# mirroring the "test 1".."test 5" rows added to resultSimulation in the
# fix described by the commit message.
$ws.Range("A2:G2").Value = "test 1"
$ws.Range("A3:G3").Value = "test 2"
$ws.Range("A4:G4").Value = "test 3"
$ws.Range("A5:G5").Value = "test 4"
$ws.Range("A6:G6").Value = "test 5"

# Match the saved selection/cursor position recorded in the workbook.
$ws.Range("J10").Select() | Out-Null
